$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row for the extra account-statement entry ---------------
# Old layout had the "closing" (thicker-border) data row at row 18, followed
# by a gap, then the footer at rows 23/24. The new layout adds one more
# data row, so the closing-style row becomes row 19 and the footer shifts
# down to rows 24/25.

$ws.Rows("19:19").Insert()

# Carry the old row 18 (closing border style + its old content) down into
# the new row 19, then restore row 18 to the "regular" border style by
# copying row 17's formatting/content into it. Values are overwritten
# below anyway; this just fixes up the per-cell styles (s attribute).
$ws.Range("B18:J18").Copy($ws.Range("B19:J19"))
$ws.Range("B17:J17").Copy($ws.Range("B18:J18"))

# --- Update the summary block ----------------------------------------------
$ws.Range("E11").Value = 107735
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 4

# --- Rewrite the data table (Tipo Doc / N Doc / Nombre / Periodo / Valor Mora / Salario) ---
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1048461061"
$ws.Range("D16").Value = "MIGUEL ANGEL LOPEZ GOMEZ"
$ws.Range("E16").Value = "2202"
$ws.Range("F16").Value = 35112
$ws.Range("G16").Value = 877803

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1048461061"
$ws.Range("D17").Value = "MIGUEL ANGEL LOPEZ GOMEZ"
$ws.Range("E17").Value = "2201"
$ws.Range("F17").Value = 35112
$ws.Range("G17").Value = 877803

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "45694727"
$ws.Range("D18").Value = "EDILMA DEL CARMEN PUCHE IZQUIERDO"
$ws.Range("E18").Value = "2009"
$ws.Range("F18").Value = 1170
$ws.Range("G18").Value = 877803

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "45689216"
$ws.Range("D19").Value = "LUZ DARY ARRIETA TINOCO"
$ws.Range("E19").Value = "2110"
$ws.Range("F19").Value = 36341
$ws.Range("G19").Value = 908526
